$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 401.72223
$ws.Range("I19").Value = 369.57144
$ws.Range("K19").Value = 369.57144
$ws.Range("M19").Value = -194.57144
$ws.Range("H55").Value = 279.41666
$ws.Range("I55").Value = 314
$ws.Range("J55").Value = 254.71428
$ws.Range("K55").Value = 314
$ws.Range("L55").Value = 254.71428
$ws.Range("M55").Value = -100
$ws.Range("N55").Value = -682.71428
$ws.Range("H112").Value = 2658.48
$ws.Range("I112").Value = 749.8333
$ws.Range("J112").Value = 3261.2104
$ws.Range("K112").Value = 2249.4999
$ws.Range("L112").Value = 9783.6312
$ws.Range("M112").Value = -1141.4999
$ws.Range("N112").Value = -11999.6312
$ws.Range("H116").Value = 2654.1304
$ws.Range("I116").Value = 2061.4167
$ws.Range("K116").Value = 2061.4167
$ws.Range("M116").Value = 1380.5833
$ws.Range("H129").Value = 838.11365
$ws.Range("J129").Value = 1017.0303
$ws.Range("L129").Value = 3051.0909
$ws.Range("N129").Value = -13051.0909
$ws.Range("H132").Value = 8136508
$ws.Range("I132").Value = 13340019
$ws.Range("K132").Value = 40020057
$ws.Range("M132").Value = -40017527
$ws.Range("H137").Value = 1078.418
$ws.Range("I137").Value = 852.7059
$ws.Range("J137").Value = 1310.9697
$ws.Range("K137").Value = 2558.1177
$ws.Range("L137").Value = 3932.9091
$ws.Range("M137").Value = -8.117700000000241
$ws.Range("N137").Value = -9032.909100000001
$ws.Range("H138").Value = 1394.48
$ws.Range("J138").Value = 1529.9493
$ws.Range("L138").Value = 4589.8479
$ws.Range("N138").Value = -14869.8479

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7578.8667
$ws.Range("I2").Value = 806.4545000000001
$ws.Range("J2").Value = 26203
$ws.Range("K2").Value = 806.4545000000001
$ws.Range("L2").Value = 26203
$ws.Range("M2").Value = -693.4545000000001
$ws.Range("N2").Value = -26429
$ws.Range("H110").Value = 1311.3182
$ws.Range("I110").Value = 1045.7693
$ws.Range("J110").Value = 1694.8889
$ws.Range("K110").Value = 1045.7693
$ws.Range("L110").Value = 1694.8889
$ws.Range("M110").Value = 999.2307000000001
$ws.Range("N110").Value = -5784.8889
$ws.Range("H116").Value = 7578.8667
$ws.Range("I116").Value = 806.4545000000001
$ws.Range("J116").Value = 26203
$ws.Range("K116").Value = 806.4545000000001
$ws.Range("L116").Value = 26203
$ws.Range("M116").Value = 1487.5455
$ws.Range("N116").Value = -30791

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7578.8667
$ws.Range("I3").Value = 806.4545000000001
$ws.Range("J3").Value = 26203
$ws.Range("K3").Value = 806.4545000000001
$ws.Range("L3").Value = 26203
$ws.Range("M3").Value = -692.4545000000001
$ws.Range("N3").Value = -26431
$ws.Range("H80").Value = 1623
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1623
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1623
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3619
$ws.Range("H83").Value = 1623
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1623
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 8115
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -18099
$ws.Range("H86").Value = 2851.1936
$ws.Range("I86").Value = 3520.6
$ws.Range("J86").Value = 1634.091
$ws.Range("K86").Value = 3520.6
$ws.Range("L86").Value = 1634.091
$ws.Range("M86").Value = -2397.6
$ws.Range("N86").Value = -3880.091
$ws.Range("H89").Value = 2851.1936
$ws.Range("I89").Value = 3520.6
$ws.Range("J89").Value = 1634.091
$ws.Range("K89").Value = 17603
$ws.Range("L89").Value = 8170.455
$ws.Range("M89").Value = -11987
$ws.Range("N89").Value = -19402.455
$ws.Range("H107").Value = 1409.9474
$ws.Range("I107").Value = 1272.6666
$ws.Range("J107").Value = 1924.75
$ws.Range("K107").Value = 1272.6666
$ws.Range("L107").Value = 1924.75
$ws.Range("M107").Value = 647.3334
$ws.Range("N107").Value = -5764.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 50395.645
$ws.Range("I22").Value = 432.7143
$ws.Range("K22").Value = 432.7143
$ws.Range("M22").Value = -82.71429999999998
$ws.Range("H31").Value = 1923.6207
$ws.Range("I31").Value = 2025.4348
$ws.Range("J31").Value = 1533.3334
$ws.Range("K31").Value = 2025.4348
$ws.Range("L31").Value = 1533.3334
$ws.Range("M31").Value = -1730.4348
$ws.Range("N31").Value = -2123.3334
$ws.Range("H34").Value = 1923.6207
$ws.Range("I34").Value = 2025.4348
$ws.Range("J34").Value = 1533.3334
$ws.Range("K34").Value = 2025.4348
$ws.Range("L34").Value = 1533.3334
$ws.Range("M34").Value = -1823.4348
$ws.Range("N34").Value = -1937.3334
$ws.Range("H58").Value = 661.58826
$ws.Range("I58").Value = 624.91174
$ws.Range("K58").Value = 624.91174
$ws.Range("M58").Value = -421.91174
$ws.Range("H99").Value = 1986.5
$ws.Range("I99").Value = 2009.25
$ws.Range("K99").Value = 2009.25
$ws.Range("M99").Value = -511.25
$ws.Range("H107").Value = 828.6667
$ws.Range("I107").Value = 787.2
$ws.Range("J107").Value = 858.2857
$ws.Range("K107").Value = 787.2
$ws.Range("L107").Value = 858.2857
$ws.Range("M107").Value = 1132.8
$ws.Range("N107").Value = -4698.2857
$ws.Range("H126").Value = 1986.5
$ws.Range("I126").Value = 2009.25
$ws.Range("K126").Value = 6027.75
$ws.Range("M126").Value = -3557.75
$ws.Range("H135").Value = 31999.166
$ws.Range("J135").Value = 31999.166
$ws.Range("L135").Value = 31999.166
$ws.Range("N135").Value = -42139.166
$ws.Range("H136").Value = 661.58826
$ws.Range("I136").Value = 624.91174
$ws.Range("K136").Value = 1874.73522
$ws.Range("M136").Value = 675.26478

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2609.625
$ws.Range("J39").Value = 2334.923
$ws.Range("L39").Value = 7004.768999999999
$ws.Range("N39").Value = -7592.768999999999
$ws.Range("H114").Value = 464.7
$ws.Range("I114").Value = 311.5
$ws.Range("J114").Value = 694.5
$ws.Range("K114").Value = 934.5
$ws.Range("L114").Value = 2083.5
$ws.Range("M114").Value = 2319.5
$ws.Range("N114").Value = -8591.5
$ws.Range("H131").Value = 23257066
$ws.Range("J131").Value = 1693.0333
$ws.Range("L131").Value = 5079.0999
$ws.Range("N131").Value = -15159.0999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56253164
$ws.Range("I70").Value = 50003324
$ws.Range("J70").Value = 66669570
$ws.Range("K70").Value = 50003324
$ws.Range("L70").Value = 66669570
$ws.Range("M70").Value = -50003054
$ws.Range("N70").Value = -66670110
$ws.Range("H73").Value = 56253164
$ws.Range("I73").Value = 50003324
$ws.Range("J73").Value = 66669570
$ws.Range("K73").Value = 50003324
$ws.Range("L73").Value = 66669570
$ws.Range("M73").Value = -50002388
$ws.Range("N73").Value = -66671442
$ws.Range("H113").Value = 1599.7142
$ws.Range("I113").Value = 1545.5
$ws.Range("K113").Value = 1545.5
$ws.Range("M113").Value = 624.5
$ws.Range("H133").Value = 50399.5
$ws.Range("J133").Value = 50399.5
$ws.Range("L133").Value = 50399.5
$ws.Range("N133").Value = -60519.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2752.5
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 2836.6667
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 2836.6667
$ws.Range("M61").Value = -2298
$ws.Range("N61").Value = -3240.6667
$ws.Range("H93").Value = 611.7
$ws.Range("I93").Value = 577.25
$ws.Range("J93").Value = 749.5
$ws.Range("K93").Value = 577.25
$ws.Range("L93").Value = 749.5
$ws.Range("M93").Value = 670.75
$ws.Range("N93").Value = -3245.5
$ws.Range("H113").Value = 2752.5
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 2836.6667
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 2836.6667
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -7176.6667
$ws.Range("H133").Value = 34466.5
$ws.Range("J133").Value = 34466.5
$ws.Range("L133").Value = 34466.5
$ws.Range("N133").Value = -39526.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3397
$ws.Range("I96").Value = 3571.4285
$ws.Range("K96").Value = 3571.4285
$ws.Range("M96").Value = -2198.4285
$ws.Range("H100").Value = 1178.75
$ws.Range("I100").Value = 1722.75
$ws.Range("J100").Value = 634.75
$ws.Range("K100").Value = 3445.5
$ws.Range("L100").Value = 1269.5
$ws.Range("M100").Value = -2904.5
$ws.Range("N100").Value = -2351.5
$ws.Range("H107").Value = 417.4
$ws.Range("I107").Value = 386
$ws.Range("K107").Value = 1158
$ws.Range("M107").Value = 762
$ws.Range("H113").Value = 559.2222
$ws.Range("I113").Value = 405.5
$ws.Range("J113").Value = 866.6667
$ws.Range("K113").Value = 1216.5
$ws.Range("L113").Value = 2600.0001
$ws.Range("M113").Value = 953.5
$ws.Range("N113").Value = -6940.0001
$ws.Range("H136").Value = 493.32
$ws.Range("I136").Value = 239.2
$ws.Range("J136").Value = 874.5
$ws.Range("K136").Value = 717.5999999999999
$ws.Range("L136").Value = 2623.5
$ws.Range("M136").Value = 1832.4
$ws.Range("N136").Value = -7723.5
